$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "pants for men sport"
$ws.Cells.Item(2, 1).Value = "mens pouch leggings"
$ws.Cells.Item(3, 1).Value = "spandex leggings boys"
$ws.Cells.Item(4, 1).Value = "rodilleras de basketball"
$ws.Cells.Item(5, 1).Value = "knee sleeves wrestling"
$ws.Cells.Item(6, 1).Value = "knee pads workout"
$ws.Cells.Item(7, 1).Value = "baseball catcher leg guards adult"
$ws.Cells.Item(8, 1).Value = "5 inch seam shorts men"
$ws.Cells.Item(9, 1).Value = "youth softball pants for girls"
$ws.Cells.Item(10, 1).Value = "mens compression workout pants"
$ws.Cells.Item(11, 1).Value = "compression shorts for men"
$ws.Cells.Item(12, 1).Value = "work knee pads for men"
$ws.Cells.Item(13, 1).Value = "padded shorts snowboarding"
$ws.Cells.Item(14, 1).Value = "youth girls yoga pants"
$ws.Cells.Item(15, 1).Value = "youth hockey padded shorts"
$ws.Cells.Item(16, 1).Value = "knee protector sports"
$ws.Cells.Item(17, 1).Value = "kneeling pad construction"
$ws.Cells.Item(18, 1).Value = "knee pad for work"
$ws.Cells.Item(19, 1).Value = "basketballs under"
$ws.Cells.Item(20, 1).Value = "boys baseball pants size 6"
$ws.Cells.Item(21, 1).Value = "sliding short"
$ws.Cells.Item(22, 1).Value = "black knee pads volleyball girls"
$ws.Cells.Item(23, 1).Value = "mens 3/4 pants"
$ws.Cells.Item(24, 1).Value = "football leg sleeves for men"
$ws.Cells.Item(25, 1).Value = "black basketball shorts men"
$ws.Cells.Item(26, 1).Value = "compression football shorts"
$ws.Cells.Item(27, 1).Value = "running leggings men"
$ws.Cells.Item(28, 1).Value = "football leggings boys"
$ws.Cells.Item(29, 1).Value = "baseball leg guards"
$ws.Cells.Item(30, 1).Value = "protective basketball"
$ws.Cells.Item(31, 1).Value = "mountain bike knee pads"
$ws.Cells.Item(32, 1).Value = "weightlifting shorts men"
$ws.Cells.Item(33, 1).Value = "shorts for men basketball"
$ws.Cells.Item(34, 1).Value = "wrestling shorts for boys"
$ws.Cells.Item(35, 1).Value = "baseball items for men"
$ws.Cells.Item(36, 1).Value = "knee pads for biking"
$ws.Cells.Item(37, 1).Value = "eva foam knee pads"
$ws.Cells.Item(38, 1).Value = "mens compression running tights"
$ws.Cells.Item(39, 1).Value = "womens softball pants black"
$ws.Cells.Item(40, 1).Value = "waist guard"
$ws.Cells.Item(41, 1).Value = "bump pads"
$ws.Cells.Item(42, 1).Value = "mens fitness pants"
$ws.Cells.Item(43, 1).Value = "cycling pants for men padded"
$ws.Cells.Item(44, 1).Value = "knee pads work"
$ws.Cells.Item(45, 1).Value = "youth football girdle"
$ws.Cells.Item(46, 1).Value = "bjj knee sleeves"
$ws.Cells.Item(47, 1).Value = "volleyball kneepads black"
$ws.Cells.Item(48, 1).Value = "non slip knee pads"
$ws.Cells.Item(49, 1).Value = "short baseball"
$ws.Cells.Item(50, 1).Value = "knee sleeves for wrestling"
$ws.Cells.Item(51, 1).Value = "knee pads"
$ws.Cells.Item(52, 1).Value = "compression sleeve youth baseball"
$ws.Cells.Item(53, 1).Value = "tights for football"
$ws.Cells.Item(54, 1).Value = "soccer pad"
$ws.Cells.Item(55, 1).Value = "cycling pants for men"
$ws.Cells.Item(56, 1).Value = "dry fit leggings men"
$ws.Cells.Item(57, 1).Value = "compression calf sleeve men basketball"
$ws.Cells.Item(58, 1).Value = "long compression shorts men"
$ws.Cells.Item(59, 1).Value = "compression shorts long men"
$ws.Cells.Item(60, 1).Value = "compression pants and tops for men"
$ws.Cells.Item(61, 1).Value = "leggings knee length"
$ws.Cells.Item(62, 1).Value = "mens softball gear"
$ws.Cells.Item(63, 1).Value = "yoga after knee replacement"
$ws.Cells.Item(64, 1).Value = "wrestling knee sleeve youth"
$ws.Cells.Item(65, 1).Value = "tights compression"
$ws.Cells.Item(66, 1).Value = "mens compression pants pack"
$ws.Cells.Item(67, 1).Value = "boys running pants"
$ws.Cells.Item(68, 1).Value = "knee pads thigh support"
$ws.Cells.Item(69, 1).Value = "youth baseball pants black"
$ws.Cells.Item(70, 1).Value = "knee pads biking adult"
$ws.Cells.Item(71, 1).Value = "youth boys leggings"
$ws.Cells.Item(72, 1).Value = "adult pants"
$ws.Cells.Item(73, 1).Value = "youth baseball compression sleeves"
$ws.Cells.Item(74, 1).Value = "calf sleeves for men football"
$ws.Cells.Item(75, 1).Value = "padded knee sleeve"
$ws.Cells.Item(76, 1).Value = "knee pad exercise"
$ws.Cells.Item(77, 1).Value = "recovery pants men"
$ws.Cells.Item(78, 1).Value = "mens tight"
$ws.Cells.Item(79, 1).Value = "mens outdoor basketball"
$ws.Cells.Item(80, 1).Value = "soccer pants youth"
$ws.Cells.Item(81, 1).Value = "protective shorts"
$ws.Cells.Item(82, 1).Value = "baseball hand guard"
$ws.Cells.Item(83, 1).Value = "bee pants"
$ws.Cells.Item(84, 1).Value = "mens protective pads"
$ws.Cells.Item(85, 1).Value = "mens shorts above knee"
$ws.Cells.Item(86, 1).Value = "basketball knee support for men"
$ws.Cells.Item(87, 1).Value = "compression pads for surgery"
$ws.Cells.Item(88, 1).Value = "snowboarding pants boys"
$ws.Cells.Item(89, 1).Value = "basketball pants for girls"
$ws.Cells.Item(90, 1).Value = "youth baseball pants girls"
$ws.Cells.Item(91, 1).Value = "knee pads for exercise"
$ws.Cells.Item(92, 1).Value = "girls compression knee sleeves"
$ws.Cells.Item(93, 1).Value = "men yoga pant"
$ws.Cells.Item(94, 1).Value = "yoga pants mens black"
$ws.Cells.Item(95, 1).Value = "basketball aids"
$ws.Cells.Item(96, 1).Value = "knee compression sleeve volleyball"
$ws.Cells.Item(97, 1).Value = "impact advanced recovery"
$ws.Cells.Item(98, 1).Value = "best basketball"
$ws.Cells.Item(99, 1).Value = "baseball compression sleeve youth"
$ws.Cells.Item(100, 1).Value = "compression knee sleeve padded"
